# Adds a new "HUCN-018 Recordatorio para los pagos" user-story block
# right after the existing "Para: Para tener recordatorio de las
# mensualidades..." paragraph (end of the HUCN-017 block), separated
# from it by one blank paragraph - mirroring the blank-line separators
# already used between the other HUCN-0xx blocks in this document.

$d = $word.ActiveDocument

# Locate the last paragraph of the HUCN-017 block (unique text).
$anchor = $d.Content
$anchor.Find.Execute(
    "Para: Para tener recordatorio de las mensualidades que lleva desde que comenzó.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $anchor.Find.Found) {
    throw "Anchor paragraph not found"
}

$insertPos = $anchor.End

# Build a zero-length range right after the anchor paragraph's text
# (NOT via Collapse(0), which in this runtime ends up swallowing the
# preceding paragraph mark) and inject the new paragraphs as raw OOXML
# so the run-splitting in "HUCN-0"/"18"/" Recordatorio..." is preserved
# exactly instead of being auto-coalesced into one run.
$insertionPoint = $d.Range($insertPos, $insertPos)

$newPartXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>HUCN-0</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>18</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> Recordatorio para los pagos</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Como: Cliente olvidadizo.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Quiero: Tener un recordatorio antes de la fecha de vencimiento para hacer el pago.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Para: Para no atrasarse en los pagos y no tener recargos extras.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint.InsertXML($newPartXml)

Write-Output "Inserted HUCN-018 block."
